# "Added last minute updates"
#
# The first paragraph of the document holds the topic-ID placeholder.
# This edit:
#   1. Updates the placeholder text to the new topic ID and drops the
#      trailing " " run that used to follow it (merged into one run).
#   2. Bumps the paragraph's left indent from 120 -> 225 twips (6pt -> 11.25pt).
#   3. Adds an (invisible, line-less) paragraph border that only carries
#      5-twip spacing on all four sides - matching the border already
#      present a couple of paragraphs further down in the document.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Replace the old placeholder (plus the trailing space that lived in the
# second run) with the new placeholder text - leaves a single run behind.
$r = $p.Range
$r.Find.Execute("**ID__AFFARS_mp_5301_602_2_d_topic_2__ID** ", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "**ID__AFFARS_MP5301_602_2_1__ID**", 2)

# Indentation: 225 twips = 11.25 points.
$p.Range.ParagraphFormat.LeftIndent = 11.25

# Paragraph border with no visible line, just 5-twip spacing on every side.
$p.Range.Borders.DistanceFromTop = 5
$p.Range.Borders.DistanceFromLeft = 5
$p.Range.Borders.DistanceFromBottom = 5
$p.Range.Borders.DistanceFromRight = 5
